$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8574223518371582
$ws.Range("B1").Value = 1.868358016014099
$ws.Range("C1").Value = 3.365347623825073
$ws.Range("D1").Value = 3.781039237976074
$ws.Range("E1").Value = 0.9734551906585693
